# disability_prevalence.xlsx (Terjola) update
# - Retitles the sheet header
# - Splits the old single "Number of disability persons" row into two rows:
#     "disabilities Persons" and "family with disabilities Persons" with new data
# - Moves the Source note down one row
# - Adjusts layout (column width, row heights, merges) to match

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row before the old row 5 (the merged "Source" row) so we end
#    up with two data rows (4 and 5) and the source note becomes row 6.
# ---------------------------------------------------------------------------
$ws.Rows("5:5").Insert()

# ---------------------------------------------------------------------------
# 2. Row 1 - title / header (merged across A1:I1)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Terjola Municipality"
$ws.Range("A1:I1").Merge()
$ws.Range("A1:I1").HorizontalAlignment = -4108
$ws.Range("A1:I1").VerticalAlignment = -4108
$ws.Range("A1:I1").WrapText = $true
$ws.Rows("1:1").RowHeight = 51

# ---------------------------------------------------------------------------
# 3. Row 2 - "(End of year, persons)" caption - content unchanged
# ---------------------------------------------------------------------------
$ws.Rows("2:2").RowHeight = 14.5

# ---------------------------------------------------------------------------
# 4. Row 3 - year headers; only A3's font changes to Sylfaen 11
# ---------------------------------------------------------------------------
$ws.Range("A3").Font.Name = "Sylfaen"
$ws.Range("A3").Font.Size = 11

# ---------------------------------------------------------------------------
# 5. Row 4 - "disabilities Persons" with new figures
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "disabilities Persons "
$ws.Range("A4").Borders.Item(9).LineStyle = -4142
$ws.Rows("4:4").RowHeight = 24.75

$ws.Range("B4").Value = 802
$ws.Range("C4").Value = 765
$ws.Range("D4").Value = 672
$ws.Range("E4").Value = 663
$ws.Range("F4").Value = 665
$ws.Range("G4").Value = 679
$ws.Range("H4").Value = 683
$ws.Range("I4").Value = 698

$ws.Range("B4:I4").NumberFormat = "#\ ##0"
$ws.Range("B4:I4").HorizontalAlignment = -4142
$ws.Range("B4:I4").Borders.Item(8).LineStyle = -4142
$ws.Range("B4:I4").Borders.Item(9).LineStyle = -4142

# ---------------------------------------------------------------------------
# 6. Row 5 (new) - "family with disabilities Persons" with new figures
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "family with disabilities Persons "
$ws.Range("A5").Font.Name = "Arial"
$ws.Range("A5").Font.Size = 10
$ws.Range("A5").Font.Bold = $false
$ws.Range("A5").HorizontalAlignment = -4131
$ws.Range("A5").VerticalAlignment = -4108
$ws.Range("A5").WrapText = $true
$ws.Range("A5").Borders.Item(9).LineStyle = 1
$ws.Range("A5").Borders.Item(9).Weight = 2
$ws.Rows("5:5").RowHeight = 21

$ws.Range("B5").Value = 935
$ws.Range("C5").Value = 889
$ws.Range("D5").Value = 804
$ws.Range("E5").Value = 797
$ws.Range("F5").Value = 809
$ws.Range("G5").Value = 816
$ws.Range("H5").Value = 821
$ws.Range("I5").Value = 833

$ws.Range("B5:I5").NumberFormat = "#\ ##0"
$ws.Range("B5:I5").HorizontalAlignment = -4142
$ws.Range("I5").Borders.Item(9).LineStyle = 1
$ws.Range("I5").Borders.Item(9).Weight = 2

# ---------------------------------------------------------------------------
# 7. Row 6 (was row 5) - Source note, now merged A6:H6, top border removed
# ---------------------------------------------------------------------------
$ws.Range("A6").Borders.Item(8).LineStyle = -4142
$ws.Rows("6:6").RowHeight = 27.75

# ---------------------------------------------------------------------------
# 8. Column / sheet level layout
# ---------------------------------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 20.8164
$ws.Range("A1").Select()
